$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "26.903.24"
$c.Style = $s

$c = $ws.Range("E2")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.76%  "
$c.Style = $s

$c = $ws.Range("D3")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.844.02"
$c.Style = $s

$c = $ws.Range("E3")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.84%  "
$c.Style = $s

$c = $ws.Range("E4")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.44%  "
$c.Style = $s

$c = $ws.Range("D5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "309.19"
$c.Style = $s

$c = $ws.Range("E5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.26%  "
$c.Style = $s

$c = $ws.Range("E6")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.31%  "
$c.Style = $s

$c = $ws.Range("D7")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.4682"
$c.Style = $s

$c = $ws.Range("E7")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +3.83%  "
$c.Style = $s

$c = $ws.Range("D8")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.3683"
$c.Style = $s

$c = $ws.Range("E8")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.86%  "
$c.Style = $s

$c = $ws.Range("D9")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.07145"
$c.Style = $s

$c = $ws.Range("E9")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.27%  "
$c.Style = $s

$c = $ws.Range("D10")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.9291"
$c.Style = $s

$c = $ws.Range("E10")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +4.55%  "
$c.Style = $s

$c = $ws.Range("D11")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "19.58"
$c.Style = $s

$c = $ws.Range("E11")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.45%  "
$c.Style = $s

$c = $ws.Range("D12")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.07695"
$c.Style = $s

$c = $ws.Range("E12")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -1.03%  "
$c.Style = $s

$c = $ws.Range("D13")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.874.45"
$c.Style = $s

$c = $ws.Range("E13")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +4.11%  "
$c.Style = $s

$c = $ws.Range("D14")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.287"
$c.Style = $s

$c = $ws.Range("E14")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.36%  "
$c.Style = $s

$c = $ws.Range("D15")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.394"
$c.Style = $s

$c = $ws.Range("E15")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.56%  "
$c.Style = $s

$c = $ws.Range("D16")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "88.10"
$c.Style = $s

$c = $ws.Range("E16")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +3.80%  "
$c.Style = $s

$c = $ws.Range("D17")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.009"
$c.Style = $s

$c = $ws.Range("E17")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.37%  "
$c.Style = $s

$c = $ws.Range("D18")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.000008629"
$c.Style = $s

$c = $ws.Range("E18")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.34%  "
$c.Style = $s

$c = $ws.Range("E19")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.25%  "
$c.Style = $s

$c = $ws.Range("D20")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "26.929.86"
$c.Style = $s

$c = $ws.Range("E20")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.68%  "
$c.Style = $s

$c = $ws.Range("D21")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "14.37"
$c.Style = $s

$c = $ws.Range("E21")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.69%  "
$c.Style = $s

$c = $ws.Range("D22")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.023"
$c.Style = $s

$c = $ws.Range("E22")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.53%  "
$c.Style = $s

$c = $ws.Range("E23")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.33%  "
$c.Style = $s

$c = $ws.Range("D24")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.930"
$c.Style = $s

$c = $ws.Range("E24")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -1.38%  "
$c.Style = $s

$c = $ws.Range("D25")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "152.38"
$c.Style = $s

$c = $ws.Range("E25")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.83%  "
$c.Style = $s

$c = $ws.Range("D26")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "18.24"
$c.Style = $s

$c = $ws.Range("E26")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.68%  "
$c.Style = $s

$c = $ws.Range("D27")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.023"
$c.Style = $s

$c = $ws.Range("E27")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -0.99%  "
$c.Style = $s

$c = $ws.Range("E28")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.75%  "
$c.Style = $s

$c = $ws.Range("D29")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.886"
$c.Style = $s

$c = $ws.Range("E29")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.41%  "
$c.Style = $s

$c = $ws.Range("D30")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.08857"
$c.Style = $s

$c = $ws.Range("E30")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.02%  "
$c.Style = $s

$c = $ws.Range("D31")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.205"
$c.Style = $s

$c = $ws.Range("E31")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.23%  "
$c.Style = $s

$c = $ws.Range("D32")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.810"
$c.Style = $s

$c = $ws.Range("E32")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.77%  "
$c.Style = $s

$c = $ws.Range("D33")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.172"
$c.Style = $s

$c = $ws.Range("E33")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +5.82%  "
$c.Style = $s

$c = $ws.Range("D34")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.7466"
$c.Style = $s

$c = $ws.Range("E34")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.23%  "
$c.Style = $s

$c = $ws.Range("D35")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.463"
$c.Style = $s

$c = $ws.Range("E35")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.93%  "
$c.Style = $s

$c = $ws.Range("D36")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.083"
$c.Style = $s

$c = $ws.Range("E36")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.64%  "
$c.Style = $s

$c = $ws.Range("D37")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.01944"
$c.Style = $s

$c = $ws.Range("E37")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.96%  "
$c.Style = $s

$c = $ws.Range("D38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.965"
$c.Style = $s

$c = $ws.Range("E38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.30%  "
$c.Style = $s

$c = $ws.Range("D39")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.05188"
$c.Style = $s

$c = $ws.Range("E39")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.10%  "
$c.Style = $s

$c = $ws.Range("D40")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.5204"
$c.Style = $s

$c = $ws.Range("D41")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.896"
$c.Style = $s

$c = $ws.Range("E41")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.42%  "
$c.Style = $s

$c = $ws.Range("D42")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.1518"
$c.Style = $s

$c = $ws.Range("E42")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.94%  "
$c.Style = $s

$c = $ws.Range("D43")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "8.130"
$c.Style = $s

$c = $ws.Range("E43")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +1.19%  "
$c.Style = $s

$c = $ws.Range("D44")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "10.53"
$c.Style = $s

$c = $ws.Range("E44")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +5.92%  "
$c.Style = $s

$c = $ws.Range("D45")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.4700"
$c.Style = $s

$c = $ws.Range("E45")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  -0.24%  "
$c.Style = $s

$c = $ws.Range("D46")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.Style = $s

$c = $ws.Range("E46")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.40%  "
$c.Style = $s

$c = $ws.Range("D47")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "100.67"
$c.Style = $s

$c = $ws.Range("E47")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.65%  "
$c.Style = $s

$c = $ws.Range("D48")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.606"
$c.Style = $s

$c = $ws.Range("E48")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.12%  "
$c.Style = $s

$c = $ws.Range("D49")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "64.82"
$c.Style = $s

$c = $ws.Range("E49")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +2.19%  "
$c.Style = $s

$c = $ws.Range("D50")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.06029"
$c.Style = $s

$c = $ws.Range("E50")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +0.78%  "
$c.Style = $s

$c = $ws.Range("D51")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.8904"
$c.Style = $s

$c = $ws.Range("E51")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "  +5.52%  "
$c.Style = $s

